$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column cells we are about to touch as Text so that
# numeric-looking strings (e.g. "0.9955", "31.229.91") are preserved exactly
# as typed instead of being parsed into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '31.229.91'
$ws.Range("E2").Value = '  +2.84%  '

$ws.Range("D3").Value = '1.979.71'
$ws.Range("E3").Value = '  +5.83%  '

$ws.Range("D4").Value = '0.9955'
$ws.Range("E4").Value = '  -0.46%  '

$ws.Range("D5").Value = '0.8170'
$ws.Range("E5").Value = '  +73.62%  '

$ws.Range("D6").Value = '253.10'
$ws.Range("E6").Value = '  +3.87%  '

$ws.Range("D7").Value = '0.9954'
$ws.Range("E7").Value = '  -0.47%  '

$ws.Range("D8").Value = '0.3407'
$ws.Range("E8").Value = '  +18.46%  '

$ws.Range("D9").Value = '25.70'
$ws.Range("E9").Value = '  +16.93%  '

$ws.Range("D10").Value = '0.06948'
$ws.Range("E10").Value = '  +7.69%  '

$ws.Range("D11").Value = '0.8383'
$ws.Range("E11").Value = '  +15.84%  '

$ws.Range("D12").Value = '0.08096'
$ws.Range("E12").Value = '  +4.23%  '

$ws.Range("D13").Value = '1.985.79'
$ws.Range("E13").Value = '  +6.36%  '

$ws.Range("D14").Value = '99.90'
$ws.Range("E14").Value = '  +4.07%  '

$ws.Range("D15").Value = '5.493'
$ws.Range("E15").Value = '  +7.36%  '

$ws.Range("D16").Value = '271.91'
$ws.Range("E16").Value = '  -2.54%  '

$ws.Range("D17").Value = '31.219.06'
$ws.Range("E17").Value = '  +2.84%  '

$ws.Range("D18").Value = '13.91'
$ws.Range("E18").Value = '  +7.17%  '

$ws.Range("D19").Value = '0.000007948'
$ws.Range("E19").Value = '  +5.84%  '

$ws.Range("D20").Value = '5.768'
$ws.Range("E20").Value = '  +10.46%  '

$ws.Range("D21").Value = '2.246.08'
$ws.Range("E21").Value = '  +6.37%  '

$ws.Range("D22").Value = '0.9965'
$ws.Range("E22").Value = '  -0.34%  '

$ws.Range("D23").Value = '0.9950'
$ws.Range("E23").Value = '  -0.51%  '

$ws.Range("D24").Value = '6.924'
$ws.Range("E24").Value = '  +11.31%  '

$ws.Range("D25").Value = '9.711'
$ws.Range("E25").Value = '  +7.42%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1511'
$ws.Range("E26").Value = '  +57.16%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '163.54'
$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("D28").Value = '19.86'
$ws.Range("E28").Value = '  +6.52%  '

$ws.Range("D29").Value = '2.188'
$ws.Range("E29").Value = '  +16.71%  '

$ws.Range("D30").Value = '1.563'
$ws.Range("E30").Value = '  +5.21%  '

$ws.Range("D31").Value = '4.564'
$ws.Range("E31").Value = '  +8.38%  '

$ws.Range("D32").Value = '1.346'
$ws.Range("E32").Value = '  +1.91%  '

$ws.Range("D33").Value = '4.310'
$ws.Range("E33").Value = '  +4.98%  '

$ws.Range("D34").Value = '0.05143'
$ws.Range("E34").Value = '  +7.06%  '

$ws.Range("E35").Value = '  +8.24%  '

$ws.Range("D36").Value = '0.7568'
$ws.Range("E36").Value = '  +10.07%  '

$ws.Range("E37").Value = '  +1.73%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02000'
$ws.Range("E38").Value = '  +6.63%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.904'
$ws.Range("E39").Value = '  +3.36%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.594'
$ws.Range("E40").Value = '  +5.99%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '78.12'
$ws.Range("E41").Value = '  +5.26%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.4674'
$ws.Range("E42").Value = '  +10.82%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '2.051'
$ws.Range("E43").Value = '  +6.27%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.8521'
$ws.Range("E44").Value = '  +3.03%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '104.22'
$ws.Range("E45").Value = '  +3.29%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '0.9951'
$ws.Range("E46").Value = '  -0.41%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '9.996'
$ws.Range("E47").Value = '  +4.65%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.484'
$ws.Range("E48").Value = '  +7.71%  '

$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").Value = '0.4288'
$ws.Range("E49").Value = '  +9.62%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '36.60'
$ws.Range("E50").Value = '  +3.74%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.1188'
$ws.Range("E51").Value = '  +12.03%  '

